$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 112279532
$ws.Range("B13").Value = 56476
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 102611
$ws.Range("F13").Value = 'Stenfalk'
$ws.Range("G13").Value = 'Falco columbarius'
$ws.Range("H13").Value = 'Linnaeus, 1758'
$ws.Range("I13").NumberFormat = "@"
$ws.Range("I13").Value = '1'
$ws.Range("I13").ClearFormats()
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = ""
$ws.Range("P13").Value = 'Arvnäs, Ly lm'
$ws.Range("Q13").Value = 694363
$ws.Range("R13").Value = 7166030
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Västerbotten'
$ws.Range("U13").Value = 'Lycksele'
$ws.Range("V13").Value = 'Lycksele lappmark'
$ws.Range("W13").Value = 'Lycksele'
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = '2023-09-14'
$ws.Range("Y13").ClearFormats()
$ws.Range("AA13").NumberFormat = "@"
$ws.Range("AA13").Value = '2023-09-14'
$ws.Range("AA13").ClearFormats()
$ws.Range("AD13").Value = $False
$ws.Range("AE13").Value = $False
$ws.Range("AG13").Value = $False
$ws.Range("AT13").Value = ""
$ws.Range("AW13").Value = 'Monika Berg'
$ws.Range("AX13").Value = 'Monika Berg'
$ws.Range("AY13").Value = ""

# Row 14
$ws.Range("A14").Value = 112279516
$ws.Range("B14").Value = 56321
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 100001
$ws.Range("F14").Value = 'Duvhök'
$ws.Range("G14").Value = 'Accipiter gentilis'
$ws.Range("H14").Value = '(Linnaeus, 1758)'
$ws.Range("I14").NumberFormat = "@"
$ws.Range("I14").Value = '1'
$ws.Range("I14").ClearFormats()
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("M14").Value = 'födosökande'
$ws.Range("N14").Value = ""
$ws.Range("P14").Value = 'Arvnäs, Ly lm'
$ws.Range("Q14").Value = 694363
$ws.Range("R14").Value = 7166030
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Västerbotten'
$ws.Range("U14").Value = 'Lycksele'
$ws.Range("V14").Value = 'Lycksele lappmark'
$ws.Range("W14").Value = 'Lycksele'
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = '2023-09-14'
$ws.Range("Y14").ClearFormats()
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = '2023-09-14'
$ws.Range("AA14").ClearFormats()
$ws.Range("AD14").Value = $False
$ws.Range("AE14").Value = $False
$ws.Range("AG14").Value = $False
$ws.Range("AT14").Value = ""
$ws.Range("AW14").Value = 'Monika Berg'
$ws.Range("AX14").Value = 'Monika Berg'
$ws.Range("AY14").Value = ""

# Row 15
$ws.Range("A15").Value = 112281233
$ws.Range("B15").Value = 57042
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'EN'
$ws.Range("E15").Value = 103042
$ws.Range("F15").Value = 'Grönfink'
$ws.Range("G15").Value = 'Chloris chloris'
$ws.Range("H15").Value = '(Linnaeus, 1758)'
$ws.Range("I15").Value = ""
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = 'födosökande'
$ws.Range("N15").Value = ""
$ws.Range("P15").Value = 'Arvnäs, Ly lm'
$ws.Range("Q15").Value = 694363
$ws.Range("R15").Value = 7166030
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = 'Västerbotten'
$ws.Range("U15").Value = 'Lycksele'
$ws.Range("V15").Value = 'Lycksele lappmark'
$ws.Range("W15").Value = 'Lycksele'
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = '2023-09-14'
$ws.Range("Y15").ClearFormats()
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = '2023-09-14'
$ws.Range("AA15").ClearFormats()
$ws.Range("AD15").Value = $False
$ws.Range("AE15").Value = $False
$ws.Range("AG15").Value = $False
$ws.Range("AT15").Value = ""
$ws.Range("AW15").Value = 'Monika Berg'
$ws.Range("AX15").Value = 'Monika Berg'
$ws.Range("AY15").Value = ""

# Row 16
$ws.Range("A16").Value = 112292314
$ws.Range("B16").Value = 56446
$ws.Range("C16").Value = 'Ovaliderad'
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 100049
$ws.Range("F16").Value = 'Spillkråka'
$ws.Range("G16").Value = 'Dryocopus martius'
$ws.Range("H16").Value = '(Linnaeus, 1758)'
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I16").Value = '1'
$ws.Range("I16").ClearFormats()
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = 'lockläte, övriga läten'
$ws.Range("N16").Value = ""
$ws.Range("P16").Value = 'Arvnäs, Ly lm'
$ws.Range("Q16").Value = 694363
$ws.Range("R16").Value = 7166030
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = 'Västerbotten'
$ws.Range("U16").Value = 'Lycksele'
$ws.Range("V16").Value = 'Lycksele lappmark'
$ws.Range("W16").Value = 'Lycksele'
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = '2023-09-15'
$ws.Range("Y16").ClearFormats()
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = '2023-09-15'
$ws.Range("AA16").ClearFormats()
$ws.Range("AD16").Value = $False
$ws.Range("AE16").Value = $False
$ws.Range("AG16").Value = $False
$ws.Range("AT16").Value = ""
$ws.Range("AW16").Value = 'Monika Berg'
$ws.Range("AX16").Value = 'Monika Berg'
$ws.Range("AY16").Value = ""

# Row 17
$ws.Range("A17").Value = 112279542
$ws.Range("B17").Value = 56841
$ws.Range("C17").Value = 'Ovaliderad'
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 103001
$ws.Range("F17").Value = 'Rödvingetrast'
$ws.Range("G17").Value = 'Turdus iliacus'
$ws.Range("H17").Value = 'Linnaeus, 1766'
$ws.Range("I17").Value = ""
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("P17").Value = 'Arvnäs, Ly lm'
$ws.Range("Q17").Value = 694363
$ws.Range("R17").Value = 7166030
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = 'Västerbotten'
$ws.Range("U17").Value = 'Lycksele'
$ws.Range("V17").Value = 'Lycksele lappmark'
$ws.Range("W17").Value = 'Lycksele'
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = '2023-09-14'
$ws.Range("Y17").ClearFormats()
$ws.Range("AA17").NumberFormat = "@"
$ws.Range("AA17").Value = '2023-09-14'
$ws.Range("AA17").ClearFormats()
$ws.Range("AD17").Value = $False
$ws.Range("AE17").Value = $False
$ws.Range("AG17").Value = $False
$ws.Range("AT17").Value = ""
$ws.Range("AW17").Value = 'Monika Berg'
$ws.Range("AX17").Value = 'Monika Berg'
$ws.Range("AY17").Value = ""

# Row 18
$ws.Range("A18").Value = 112281154
$ws.Range("B18").Value = 56575
$ws.Range("C18").Value = 'Ovaliderad'
$ws.Range("D18").Value = 'NT'
$ws.Range("E18").Value = 103021
$ws.Range("F18").Value = 'Talltita'
$ws.Range("G18").Value = 'Poecile montanus'
$ws.Range("H18").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I18").Value = ""
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = 'födosökande'
$ws.Range("N18").Value = ""
$ws.Range("P18").Value = 'Arvnäs, Ly lm'
$ws.Range("Q18").Value = 694363
$ws.Range("R18").Value = 7166030
$ws.Range("S18").Value = 10
$ws.Range("T18").Value = 'Västerbotten'
$ws.Range("U18").Value = 'Lycksele'
$ws.Range("V18").Value = 'Lycksele lappmark'
$ws.Range("W18").Value = 'Lycksele'
$ws.Range("Y18").NumberFormat = "@"
$ws.Range("Y18").Value = '2023-09-14'
$ws.Range("Y18").ClearFormats()
$ws.Range("AA18").NumberFormat = "@"
$ws.Range("AA18").Value = '2023-09-14'
$ws.Range("AA18").ClearFormats()
$ws.Range("AD18").Value = $False
$ws.Range("AE18").Value = $False
$ws.Range("AG18").Value = $False
$ws.Range("AT18").Value = ""
$ws.Range("AW18").Value = 'Monika Berg'
$ws.Range("AX18").Value = 'Monika Berg'
$ws.Range("AY18").Value = ""

# Row 19
$ws.Range("A19").Value = 112281199
$ws.Range("B19").Value = 57103
$ws.Range("C19").Value = 'Ovaliderad'
$ws.Range("D19").Value = 'NT'
$ws.Range("E19").Value = 103057
$ws.Range("F19").Value = 'Sävsparv'
$ws.Range("G19").Value = 'Emberiza schoeniclus'
$ws.Range("H19").Value = '(Linnaeus, 1758)'
$ws.Range("I19").Value = ""
$ws.Range("K19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = ""
$ws.Range("P19").Value = 'Arvnäs, Ly lm'
$ws.Range("Q19").Value = 694363
$ws.Range("R19").Value = 7166030
$ws.Range("S19").Value = 10
$ws.Range("T19").Value = 'Västerbotten'
$ws.Range("U19").Value = 'Lycksele'
$ws.Range("V19").Value = 'Lycksele lappmark'
$ws.Range("W19").Value = 'Lycksele'
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("Y19").Value = '2023-09-14'
$ws.Range("Y19").ClearFormats()
$ws.Range("AA19").NumberFormat = "@"
$ws.Range("AA19").Value = '2023-09-14'
$ws.Range("AA19").ClearFormats()
$ws.Range("AD19").Value = $False
$ws.Range("AE19").Value = $False
$ws.Range("AG19").Value = $False
$ws.Range("AT19").Value = ""
$ws.Range("AW19").Value = 'Monika Berg'
$ws.Range("AX19").Value = 'Monika Berg'
$ws.Range("AY19").Value = ""

# Row 20
$ws.Range("A20").Value = 112279543
$ws.Range("B20").Value = 56847
$ws.Range("C20").Value = 'Ovaliderad'
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 102999
$ws.Range("F20").Value = 'Björktrast'
$ws.Range("G20").Value = 'Turdus pilaris'
$ws.Range("H20").Value = 'Linnaeus, 1758'
$ws.Range("I20").Value = ""
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("P20").Value = 'Arvnäs, Ly lm'
$ws.Range("Q20").Value = 694363
$ws.Range("R20").Value = 7166030
$ws.Range("S20").Value = 10
$ws.Range("T20").Value = 'Västerbotten'
$ws.Range("U20").Value = 'Lycksele'
$ws.Range("V20").Value = 'Lycksele lappmark'
$ws.Range("W20").Value = 'Lycksele'
$ws.Range("Y20").NumberFormat = "@"
$ws.Range("Y20").Value = '2023-09-14'
$ws.Range("Y20").ClearFormats()
$ws.Range("AA20").NumberFormat = "@"
$ws.Range("AA20").Value = '2023-09-14'
$ws.Range("AA20").ClearFormats()
$ws.Range("AD20").Value = $False
$ws.Range("AE20").Value = $False
$ws.Range("AG20").Value = $False
$ws.Range("AT20").Value = ""
$ws.Range("AW20").Value = 'Monika Berg'
$ws.Range("AX20").Value = 'Monika Berg'
$ws.Range("AY20").Value = ""

# Row 21
$ws.Range("A21").Value = 112281210
$ws.Range("B21").Value = 57076
$ws.Range("C21").Value = 'Ovaliderad'
$ws.Range("D21").Value = 'VU'
$ws.Range("E21").Value = 103053
$ws.Range("F21").Value = 'Lappsparv'
$ws.Range("G21").Value = 'Calcarius lapponicus'
$ws.Range("H21").Value = '(Linnaeus, 1758)'
$ws.Range("I21").NumberFormat = "@"
$ws.Range("I21").Value = '1'
$ws.Range("I21").ClearFormats()
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = 'födosökande'
$ws.Range("N21").Value = ""
$ws.Range("P21").Value = 'Arvnäs, Ly lm'
$ws.Range("Q21").Value = 694363
$ws.Range("R21").Value = 7166030
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = 'Västerbotten'
$ws.Range("U21").Value = 'Lycksele'
$ws.Range("V21").Value = 'Lycksele lappmark'
$ws.Range("W21").Value = 'Lycksele'
$ws.Range("Y21").NumberFormat = "@"
$ws.Range("Y21").Value = '2023-09-14'
$ws.Range("Y21").ClearFormats()
$ws.Range("AA21").NumberFormat = "@"
$ws.Range("AA21").Value = '2023-09-14'
$ws.Range("AA21").ClearFormats()
$ws.Range("AD21").Value = $False
$ws.Range("AE21").Value = $False
$ws.Range("AG21").Value = $False
$ws.Range("AT21").Value = ""
$ws.Range("AW21").Value = 'Monika Berg'
$ws.Range("AX21").Value = 'Monika Berg'
$ws.Range("AY21").Value = ""

Write-Output "Edit complete"
